# Upload new version with timestamp
# Adds two new inventory rows (item 4 "SPASMO-DIGESTIN..." swaps into the
# existing item-3 slot, "SUGARLO PLUS ..." becomes a new item 4, and the
# former item 3 "اولويز..." is pushed down to item 5), refreshes the
# totals row and bumps the generated-at timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$val) {
    # Force literal text storage (shared string) while preserving whatever
    # number format the cell already had (mirrors cells that show numeric
    # looking text, e.g. "1", "78.00", stored as text, not numbers).
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.NumberFormat = $fmt
}

# --- Step 1: insert two fresh rows right above the totals row ------------
# (old row 10 "totals" / row 11 "footer" shift down to 12 / 13)
$ws.Rows("10:11").Insert()

# --- Step 2: clone formatting from the existing item row (row 9) ---------
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: rebuild the merged cells for the two new rows ----------------
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# --- Step 4: row heights ---------------------------------------------------
$ws.Rows(9).RowHeight = 25.5
$ws.Rows(10).RowHeight = 24.75
$ws.Rows(11).RowHeight = 25.5
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 16.5

# --- Step 5: row 9 now holds "SPASMO-DIGESTIN 30 TABS." (was item 3) ------
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "SPASMO-DIGESTIN 30 TABS."
$ws.Range("H9").Value = "4:0"
Set-TextValue $ws.Range("L9") "1"
$ws.Range("N9").Value = "78.00"
Set-TextValue $ws.Range("P9") "25.7400"
$ws.Range("Q9").Value = "0:1"

# --- Step 6: row 10 is the new item 4 "SUGARLO PLUS ..." -------------------
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "SUGARLO PLUS 50/1000MG 30 F.C. TABS"
$ws.Range("H10").Value = "0:0"
Set-TextValue $ws.Range("L10") "1"
$ws.Range("N10").Value = "136.50"
Set-TextValue $ws.Range("P10") "45.0450"
$ws.Range("Q10").Value = "0:1"

# --- Step 7: row 11 is item 5, the original item 3 "اولويز..." pushed down -
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "اولويز ماكس طويل جدا"
$ws.Range("H11").Value = "43:0"
Set-TextValue $ws.Range("L11") "0"
$ws.Range("N11").Value = "35.00"
Set-TextValue $ws.Range("P11") "35.0000"
$ws.Range("Q11").Value = "1:0"

# --- Step 8: refresh the totals row (now row 12) ---------------------------
$ws.Range("P12").Value = 160.04499999999999

# --- Step 9: bump the generated-at timestamp in the footer (now row 13) ---
$ws.Range("A13").Value = "Monday, 25 August, 2025 10:18 AM"

Write-Host "edit applied"
